# Add "bills_cards" and "teams" worksheets (with their table_column_names
# data) after the existing "cards" sheet, and update the various sheet
# selections / active tab to match the final state of the workbook.

$wb = $excel.ActiveWorkbook

$wsPlayers = $wb.Worksheets.Item("players")
$wsCards   = $wb.Worksheets.Item("cards")

# ---------------------------------------------------------------------
# 1. Create the two new sheets, in order, right after "cards".
# ---------------------------------------------------------------------
$wsBills = $wb.Worksheets.Add($null, $wsCards)
$wsBills.Name = "bills_cards"

$wsTeams = $wb.Worksheets.Add($null, $wsBills)
$wsTeams.Name = "teams"

# ---------------------------------------------------------------------
# 2. bills_cards sheet: same look as "cards" (header row + PK/U/FK style
#    column) but describing the bills_card table's columns.
# ---------------------------------------------------------------------

# Header row (COLUMN_NAME / NULL / DATA_TYPE banner, style copied from cards!A1:D1)
$wsCards.Range("A1:D1").Copy($wsBills.Range("A1:D1"))

# Body rows 2-11 reuse the "cards!A2:D2" look (col A = bold/no-fill/bordered,
# cols B:D = plain/bordered) then row 12 reuses "cards!A8:D8" (the FK row).
for ($r = 2; $r -le 11; $r++) {
    $wsCards.Range("A2:D2").Copy($wsBills.Range("A$r`:D$r"))
}
$wsCards.Range("A8:D8").Copy($wsBills.Range("A12:D12"))

# Now fill in the real values. NOTE: the row order here (3..12, then 2
# last) is deliberate - it reproduces the original author's shared-string
# insertion order (the "bills_card_id" PK row was filled in last).
$billsRows = @(
    @(3,  "U",  "bills_card_cert",      "NO",  "int"),
    @(4,  "",   "bills_card_spec",      "NO",  "int"),
    @(5,  "",   "bills_card_num",       "NO",  "int"),
    @(6,  "",   "bills_card_year",      "NO",  "int"),
    @(7,  "",   "bills_card_psa_desc",  "NO",  "nvarchar"),
    @(8,  "",   "bills_card_grade",     "NO",  "float"),
    @(9,  "",   "bills_card_pop",       "NO",  "int"),
    @(10, "",   "bills_card_pop_higher","NO",  "int"),
    @(11, "",   "bills_card_stat_year", "YES", "int"),
    @(12, "FK", "bills_card_player_id", "YES", "int"),
    @(2,  "PK", "bills_card_id",        "NO",  "int")
)

foreach ($row in $billsRows) {
    $r = $row[0]
    if ($row[1] -ne "") {
        $wsBills.Range("A$r").Value = $row[1]
    } else {
        $wsBills.Range("A$r").ClearContents()
    }
    $wsBills.Range("B$r").Value = $row[2]
    $wsBills.Range("C$r").Value = $row[3]
    $wsBills.Range("D$r").Value = $row[4]
}

$wsBills.Columns.Item(1).ColumnWidth = 3.96
$wsBills.Columns.Item(2).ColumnWidth = 18.11
$wsBills.Columns.Item(3).ColumnWidth = 6.17
$wsBills.Columns.Item(4).ColumnWidth = 9.37

# ---------------------------------------------------------------------
# 3. teams sheet: same look, describing the team table's columns.
# ---------------------------------------------------------------------

# Header row: copy format+value for B1:D1, but A1 has no PK/etc. text and
# uses the "plain header" look (bold cleared) rather than bold.
$wsCards.Range("A1:D1").Copy($wsTeams.Range("A1:D1"))
$wsTeams.Range("A1").ClearContents()
$wsTeams.Range("A1").Font.Bold = $false

# Row 2 (PK) reuses cards!A2:D2 look; rows 3-5 reuse cards!A6:D6 (plain, no
# PK/U/FK marker) look.
$wsCards.Range("A2:D2").Copy($wsTeams.Range("A2:D2"))
for ($r = 3; $r -le 5; $r++) {
    $wsCards.Range("A6:D6").Copy($wsTeams.Range("A$r`:D$r"))
}

$teamsRows = @(
    @("PK", "team_code",   "NO", "char"),
    @("",   "team_city",   "NO", "nvarchar"),
    @("",   "team_name",   "NO", "nvarchar"),
    @("",   "team_league", "NO", "nvarchar")
)

$r = 2
foreach ($row in $teamsRows) {
    if ($row[0] -ne "") {
        $wsTeams.Range("A$r").Value = $row[0]
    } else {
        $wsTeams.Range("A$r").ClearContents()
    }
    $wsTeams.Range("B$r").Value = $row[1]
    $wsTeams.Range("C$r").Value = $row[2]
    $wsTeams.Range("D$r").Value = $row[3]
    $r++
}

$wsTeams.Columns.Item(1).ColumnWidth = 3.22
$wsTeams.Columns.Item(2).ColumnWidth = 13.27
$wsTeams.Columns.Item(3).ColumnWidth = 6.01
$wsTeams.Columns.Item(4).ColumnWidth = 9.37

# ---------------------------------------------------------------------
# 4. Selections / active tab. Order matters since .Select() activates
#    the sheet it's called on - end on "teams" so it is the active tab.
# ---------------------------------------------------------------------

$wsCards.Activate() | Out-Null
$wsCards.Range("E10").Select() | Out-Null

$wsBills.Activate() | Out-Null
$wsBills.Range("A1:D12").Select() | Out-Null

$wsTeams.Activate() | Out-Null
$wsTeams.Range("A1:D5").Select() | Out-Null
